$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '24.839.14'
$ws.Range("E2").Value = '  +0.22%  '

$ws.Range("D3").Value = '1.686.91'
$ws.Range("E3").Value = '  -1.15%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.006'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.87%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '314.79'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.78%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3944'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.23%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3971'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -2.31%  '

$ws.Range("B9").Value = 'BinanceUSD'
$ws.Range("C9").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '1.006'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.09%  '

$ws.Range("B10").Value = 'Polygon'
$ws.Range("C10").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.426'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -4.91%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '51.68'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -3.74%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.08681'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -1.59%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '25.12'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -4.38%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.290'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -2.65%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.811'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -3.93%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.00001321'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -3.03%  '

$ws.Range("D17").Value = '1.605.15'
$ws.Range("E17").Value = '  -5.49%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '94.09'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -3.72%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.07129'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.40%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '20.20'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.82%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.183'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.63%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.005'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.79%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '14.12'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.80%  '

$ws.Range("D24").Value = '24.860.13'
$ws.Range("E24").Value = '  +0.34%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.395'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +2.54%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.786'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -7.64%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '23.18'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.77%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '5.960'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.31%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '161.39'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -3.70%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '149.08'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +2.64%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.639'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +21.37%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.847'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -8.88%  '

$ws.Range("D33").Value = '1.814.67'
$ws.Range("E33").Value = '  -3.81%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.08441'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -4.34%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.03075'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.33%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.015'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -5.15%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '6.948'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -3.85%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.2826'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.24%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.09580'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +3.95%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '10.59'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -3.06%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.8006'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -6.28%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '13.72'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -3.38%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.458'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.17%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '16.77'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -5.05%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.595'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -4.45%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.7192'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -4.42%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '4.199'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.94%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.08749'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +5.69%  '

$ws.Range("B49").Value = 'Frax'
$ws.Range("C49").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.003'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.66%  '

$ws.Range("B50").Value = 'Flow'
$ws.Range("C50").Value = 'https://coinranking.com/coin/QQ0NCmjVq+flow-flow'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.351'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.76%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '139.01'
$ws.Range("D51").Style = "Normal"

